# Update the "Förändrad" (changed) date column (C) for rows 2-79
# from serial date 45175 (2023-09-06) to 45177 (2023-09-08).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C79").Value = 45177
